# "Generate Report for Handback"
#
# This script updates the localization-status report to reflect that the
# handback files are now in sync with en-US (instead of just "Ready for
# handoff"), refreshes the "Latest Handback DateTime" timestamps for the
# zh-cn and de-de sheets, and clears the stale "Error Detail" message that
# said the handback file was out of date (it no longer applies once the
# handback is in sync).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de status columns (E, F) both show the
# shared "Ready for handoff" text -> update to the new status text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Widen the status columns so the longer text fits.
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 29.14

# Latest Handback DateTime refreshed for this handback report
$wsZhCn.Range("K2").Value = "2016-10-20 09:52:45"
$wsZhCn.Range("K3").Value = "2016-10-20 09:52:45"

# Error Detail is no longer applicable now that the handback is current.
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 29.14

# Latest Handback DateTime refreshed for this handback report
$wsDeDe.Range("K2").Value = "2016-10-20 09:53:03"
$wsDeDe.Range("K3").Value = "2016-10-20 09:53:03"

# Error Detail is no longer applicable now that the handback is current.
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8
